# Apply commit: update France Commemorative UNC varieties workbook
# - Rename "UNICEF" subject (row 33) to "75th Anniversary - UNICEF"
# - Bump several "duplicates" counters (I column) from 0/1 to 1/2
# - Add a note on I21 explaining the duplicate is not UNC
# - Move the active selection to I4

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("2€")

# Rename the Subject text for row 33 (2021 UNICEF coin)
$ws.Range("B33").Value = "75th Anniversary - UNICEF"

# Update duplicate counters
$ws.Range("I16").Value = 1
$ws.Range("I21").Value = 2
$ws.Range("I25").Value = 1
$ws.Range("I33").Value = 1
$ws.Range("I35").Value = 1
$ws.Range("I36").Value = 1
$ws.Range("I39").Value = 1

# Add comment to I21 explaining duplicates are not UNC
$excel.UserName = "Lord_Alexator"
$comment = $ws.Range("I21").AddComment("Дубли не UNC")

# Update the active selection on the sheet
$ws.Range("I4").Select()

$wb.Save()
